$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 3rd sheet in this workbook
$ws = $wb.Worksheets.Item(3)

# Insert a new (blank) column before the "Late" column so the schedule
# gains an extra column for the new Variable Instalments data, shifting
# the existing Late / heading / Outstanding columns one place to the right.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.14

# Make the "Repayment schedule" sheet the active tab / selection, as it
# was left selected at cell S7 after the edit.
$ws.Activate()
$ws.Range("S7").Select()
